$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 1155-1156, pushing the existing data (old rows
# 1155..1245) down to 1157..1247. Excel's row Insert mirrors formatting
# from the row above, which also gives the new D-column cells the date
# number format used throughout the column.
$ws.Rows("1155:1156").Insert()

# Fill in the two new rows with the new weekly price observations.
$ws.Range("A1155").Value = 9
$ws.Range("B1155").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C1155").Value = "Metropolitana"
$ws.Range("D1155").Value = 44746
$ws.Range("E1155").Value = 13
$ws.Range("F1155").Value = "Fruta"
$ws.Range("G1155").Value = 100108
$ws.Range("H1155").Value = "Tropicales y subtropicales"
$ws.Range("I1155").Value = 100108006
$ws.Range("J1155").Value = "Plátano"
$ws.Range("K1155").Value = "Sin especificar"
$ws.Range("L1155").Value = "Primera Maduro"
$ws.Range("M1155").Value = 630
$ws.Range("N1155").Value = 22000
$ws.Range("O1155").Value = 23000
$ws.Range("P1155").Value = 22556
$ws.Range("Q1155").Value = "$/caja 20 kilos"
$ws.Range("R1155").Value = "Ecuador"
$ws.Range("S1155").Value = 1128
$ws.Range("T1155").Value = 20

$ws.Range("A1156").Value = 9
$ws.Range("B1156").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C1156").Value = "Metropolitana"
$ws.Range("D1156").Value = 44746
$ws.Range("E1156").Value = 13
$ws.Range("F1156").Value = "Fruta"
$ws.Range("G1156").Value = 100108
$ws.Range("H1156").Value = "Tropicales y subtropicales"
$ws.Range("I1156").Value = 100108006
$ws.Range("J1156").Value = "Plátano"
$ws.Range("K1156").Value = "Sin especificar"
$ws.Range("L1156").Value = "Primera Pintón"
$ws.Range("M1156").Value = 990
$ws.Range("N1156").Value = 23000
$ws.Range("O1156").Value = 24000
$ws.Range("P1156").Value = 23455
$ws.Range("Q1156").Value = "$/caja 20 kilos"
$ws.Range("R1156").Value = "Ecuador"
$ws.Range("S1156").Value = 1173
$ws.Range("T1156").Value = 20
